$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the summary / header row (row 2) counters ---
$ws.Range("F2").Value = "31 / 77"
$ws.Range("G2").Value = "32 / 39"
$ws.Range("J2").Value = "3 / 15"

# --- Insert a new standings row at row 23, pushing rows 23:35 down to 24:36 ---
$ws.Rows.Item(23).Insert()

# Copy the formatting of the (now shifted) row 24 back onto the newly
# inserted row 23 so the new row matches the rest of the table (border,
# bold centered rank cell in column A, etc.) -- restrict to the used
# columns (A:K) so we don't balloon the sheet's dimension out to XFD.
$ws.Range("A24:K24").Copy()
$ws.Range("A23:K23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Renumber columns A (place) and B (rank) for every row that shifted down
# one position (old rows 23-35, now at 24-36): bump each by 1.
for ($r = 24; $r -le 36; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 1).Value() + 1
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 2).Value() + 1
}

# --- Fill in the brand-new row 23 with the new team's standings ---
$ws.Cells.Item(23, 1).Value = 19
$ws.Cells.Item(23, 2).Value = 20
$ws.Cells.Item(23, 3).Value = "Sherlocked_hzoi (王若竹)"
$ws.Cells.Item(23, 4).Value = 2
$ws.Cells.Item(23, 5).Value = "68  1:08:35"
$ws.Cells.Item(23, 6).Value = "0:10:29  (-2)"
$ws.Cells.Item(23, 7).Value = "0:18:06"
$ws.Cells.Item(23, 8).Value = ""
$ws.Cells.Item(23, 9).Value = ""
$ws.Cells.Item(23, 10).Value = "(-2)"
$ws.Cells.Item(23, 11).Value = ""
